$d = $word.ActiveDocument

# --- Step 1: remove the existing hidden "_GoBack" bookmark ---
# (it currently sits at the end of the "Resultaten" paragraph; it will be
#  re-created around the newly struck-through text below)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: split "Bart Muelders & Feiko Wielsma" and strike the second half ---
$find = $d.Content.Find
$find.Execute("& Feiko Wielsma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $find.Parent

$rng.Font.StrikeThrough = $true

# --- Step 3: wrap that same range with a new "_GoBack" bookmark ---
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
